$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sriram")

# Activate this sheet (it was already the tabSelected sheet)
$ws.Activate()

# Row 24 is the closest formatting template: bare A cell, date-formatted B,
# wrapped-text C/D, and the amber "in progress" E style we need for row 26.
$ws.Range("A24:E24").Copy()
$ws.Range("A26:E26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New row 20 of the log table (spreadsheet row 26)
$ws.Cells.Item(26, 1).Value = 20
$ws.Cells.Item(26, 2).Value = "2/20/2018"
$ws.Cells.Item(26, 3).Value = "11:00 to 17:00"
$ws.Cells.Item(26, 4).Value = "Worked on Receive mail and tried to filter mails with subject: IMPORTANT"
$ws.Cells.Item(26, 5).Value = "in progress"

# Row 26 is a taller "wrapped" row, like row 25 above it
$ws.Rows.Item(26).RowHeight = $ws.Rows.Item(25).RowHeight

# Update the view: scroll down, select the newly added status cell
$excel.ActiveWindow.ScrollRow = 17
$ws.Range("E26").Select()
